# Apply updated A-column sequence numbers and fix J302 value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Value=995},
    @{Row=3; Value=996},
    @{Row=4; Value=997},
    @{Row=5; Value=998},
    @{Row=6; Value=999},
    @{Row=7; Value=1000},
    @{Row=8; Value=1001},
    @{Row=11; Value=1017},
    @{Row=13; Value=1027},
    @{Row=14; Value=1028},
    @{Row=16; Value=1044},
    @{Row=17; Value=945},
    @{Row=18; Value=1050},
    @{Row=19; Value=1051},
    @{Row=21; Value=957},
    @{Row=23; Value=958},
    @{Row=24; Value=1034},
    @{Row=25; Value=1035},
    @{Row=26; Value=1036},
    @{Row=27; Value=950},
    @{Row=38; Value=1007},
    @{Row=39; Value=1008},
    @{Row=40; Value=1009},
    @{Row=41; Value=1010},
    @{Row=42; Value=1011},
    @{Row=43; Value=1012},
    @{Row=44; Value=1013},
    @{Row=45; Value=1061},
    @{Row=48; Value=1062},
    @{Row=50; Value=1018},
    @{Row=51; Value=1019},
    @{Row=52; Value=1020},
    @{Row=55; Value=1037},
    @{Row=60; Value=947},
    @{Row=64; Value=1038},
    @{Row=66; Value=1065},
    @{Row=68; Value=1014},
    @{Row=81; Value=966},
    @{Row=83; Value=967},
    @{Row=84; Value=968},
    @{Row=85; Value=969},
    @{Row=86; Value=970},
    @{Row=87; Value=971},
    @{Row=88; Value=972},
    @{Row=89; Value=973},
    @{Row=90; Value=974},
    @{Row=91; Value=975},
    @{Row=92; Value=976},
    @{Row=93; Value=977},
    @{Row=94; Value=978},
    @{Row=95; Value=979},
    @{Row=96; Value=980},
    @{Row=97; Value=981},
    @{Row=98; Value=982},
    @{Row=99; Value=983},
    @{Row=100; Value=984},
    @{Row=101; Value=985},
    @{Row=102; Value=952},
    @{Row=103; Value=1045},
    @{Row=104; Value=1046},
    @{Row=105; Value=1047},
    @{Row=111; Value=953},
    @{Row=124; Value=1048},
    @{Row=128; Value=1067},
    @{Row=133; Value=988},
    @{Row=134; Value=986},
    @{Row=138; Value=1015},
    @{Row=148; Value=1053},
    @{Row=158; Value=1016},
    @{Row=165; Value=1029},
    @{Row=166; Value=1030},
    @{Row=167; Value=1031},
    @{Row=168; Value=1032},
    @{Row=173; Value=962},
    @{Row=174; Value=943},
    @{Row=176; Value=946},
    @{Row=177; Value=987},
    @{Row=180; Value=1063},
    @{Row=181; Value=1049},
    @{Row=190; Value=991},
    @{Row=193; Value=1039},
    @{Row=194; Value=959},
    @{Row=195; Value=1002},
    @{Row=196; Value=992},
    @{Row=197; Value=993},
    @{Row=202; Value=1059},
    @{Row=203; Value=1060},
    @{Row=204; Value=1003},
    @{Row=205; Value=1004},
    @{Row=206; Value=1005},
    @{Row=207; Value=1006},
    @{Row=211; Value=1052},
    @{Row=214; Value=1064},
    @{Row=215; Value=1054},
    @{Row=216; Value=1055},
    @{Row=217; Value=1056},
    @{Row=223; Value=1021},
    @{Row=224; Value=1022},
    @{Row=231; Value=990},
    @{Row=232; Value=964},
    @{Row=233; Value=1043},
    @{Row=235; Value=954},
    @{Row=236; Value=1023},
    @{Row=237; Value=1024},
    @{Row=238; Value=1025},
    @{Row=239; Value=1026},
    @{Row=241; Value=1033},
    @{Row=243; Value=948},
    @{Row=244; Value=1040},
    @{Row=245; Value=949},
    @{Row=246; Value=944},
    @{Row=247; Value=960},
    @{Row=248; Value=951},
    @{Row=249; Value=965},
    @{Row=250; Value=1057},
    @{Row=252; Value=1066},
    @{Row=253; Value=1068},
    @{Row=259; Value=1041},
    @{Row=260; Value=989},
    @{Row=264; Value=961},
    @{Row=267; Value=963},
    @{Row=268; Value=994},
    @{Row=269; Value=1042},
    @{Row=280; Value=1058},
    @{Row=300; Value=955},
    @{Row=301; Value=956},
    @{Row=302; Value=102}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Value
}

# Fix J302 value (recallD3) from -0 to 5830
$ws.Cells.Item(302, 10).Value = 5830
